$d = $word.ActiveDocument
$cr = [char]13
$nbsp = [char]0x00A0

$oldTitle  = "SMARTCASH HIVE TEAMS"
$newTitle  = "Команды SmartCash Hive"

$oldIntro  = "SmartCash is aiming to create a decentralized team structure by efficiently distributing the workload across multiple global Hive teams. The teams below are simply getting things started. After SmartCash starts to require more teams they will be created and these will be splintered into smaller teams."
$newIntro  = "SmartCash стремится создать децентрализованную командную структуру, эффективно распределяя рабочую нагрузку между несколькими глобальными командами Hive.   Когда SmartCash потребуется большее количество команд — они будут созданы для поддержания лучшей координации и эффективности работы.`n"

$oldHeading = $nbsp + "HIVE TEAM: OUTREACH"
$newHeading = "Команда Hive: Продвижение"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $txt = $p.Range.Text.TrimEnd($cr)

    if ($txt -eq $oldTitle) {
        # Whole paragraph is a single run -> replace via the paragraph's own
        # Range so the trailing paragraph mark (and xml:space="preserve") are
        # kept intact.
        $p.Range.Text = $newTitle
    }
    elseif ($txt -eq $oldIntro) {
        $p.Range.Text = $newIntro
    }
    elseif ($txt -eq $oldHeading) {
        # This paragraph has two runs: a leading NBSP run (keep as-is) and
        # the "HIVE TEAM: OUTREACH" run. Only touch the second run so its
        # distinct formatting (and the NBSP run) survive untouched, and so
        # the sibling "...OUTREACH 2" paragraph elsewhere is never matched.
        $sub = $p.Range.Duplicate
        [void]$sub.MoveStart(1, 1)
        $sub.Text = $newHeading
    }
}
